$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text for column C
$ws.Range("C1").Value = "Change"

# Update correlation values in column C with new Yahoo Finance based data
$ws.Range("C2").Value = 0.0903
$ws.Range("C3").Value = 0.0755
$ws.Range("C4").Value = 0.0419
$ws.Range("C5").Value = 0.0979
$ws.Range("C6").Value = 0.0947
